$wb = $excel.ActiveWorkbook

# ---------- Sheet: LP1912 ----------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 07:59:05'
$ws1.Range("A3").Value = 'Total filas: 93'

$data1 = @(
    @('04:01:01', '04:01', '81_EL PELIGRO', 0, 'LP1912'),
    @('04:36:47', '04:45', '215A_EL PATO', 9, 'LP1912'),
    @('04:01:01', '04:46', '215A_EL PATO', 45, 'LP1912'),
    @('04:01:01', '04:53', '11_ETCHEVERRY', 52, 'LP1912'),
    @('04:51:28', '05:13', '14_ABASTO', 22, 'LP1912'),
    @('04:36:47', '05:14', '14_ABASTO', 38, 'LP1912'),
    @('04:01:01', '05:16', '17_ROMERO', 75, 'LP1912'),
    @('05:20:00', '05:22', '14_ABASTO', 2, 'LP1912'),
    @('04:01:01', '05:22', '23_HERNANDEZ', 81, 'LP1912'),
    @('04:36:47', '05:34', '215B_EL PATO', 58, 'LP1912'),
    @('04:01:01', '05:35', '215B_EL PATO', 94, 'LP1912'),
    @('04:01:01', '05:41', '14_ABASTO', 100, 'LP1912'),
    @('04:01:01', '05:46', '15_ABASTO', 105, 'LP1912'),
    @('04:36:47', '06:04', '16_SANTA ANA', 88, 'LP1912'),
    @('04:36:47', '06:11', '215A_EL PATO', 95, 'LP1912'),
    @('05:20:00', '06:12', '215A_EL PATO', 52, 'LP1912'),
    @('04:36:47', '06:14', '225_HARAS DEL SUR', 98, 'LP1912'),
    @('04:36:47', '06:21', '26_HERNANDEZ', 105, 'LP1912'),
    @('04:36:47', '06:27', '23_HERNANDEZ', 111, 'LP1912'),
    @('04:36:47', '06:29', '86_EST CHICA-ESC AGRARIA', 113, 'LP1912'),
    @('05:20:00', '06:30', '86_EST CHICA-ESC AGRARIA', 70, 'LP1912'),
    @('04:36:47', '06:31', '16_SANTA ANA', 115, 'LP1912'),
    @('06:23:52', '06:33', '23_HERNANDEZ', 10, 'LP1912'),
    @('06:23:52', '06:43', '225_C ROCA-H SUR', 20, 'LP1912'),
    @('04:51:28', '06:44', '225_C ROCA-H SUR', 113, 'LP1912'),
    @('04:51:28', '06:46', '215C_EL PATO', 115, 'LP1912'),
    @('05:20:00', '06:47', '215C_EL PATO', 87, 'LP1912'),
    @('06:52:23', '06:52', '215C_EL PATO', 0, 'LP1912'),
    @('05:51:32', '06:59', '14_ABASTO', 68, 'LP1912'),
    @('05:20:00', '07:00', '10_OLMOS', 100, 'LP1912'),
    @('05:20:00', '07:00', '14_ABASTO', 100, 'LP1912'),
    @('06:23:52', '07:01', '16_SANTA ANA', 38, 'LP1912'),
    @('05:51:32', '07:04', '23_HERNANDEZ', 73, 'LP1912'),
    @('05:20:00', '07:05', '15_ABASTO', 105, 'LP1912'),
    @('06:23:52', '07:06', '225_GOMEZ', 43, 'LP1912'),
    @('05:20:00', '07:07', '225_GOMEZ', 107, 'LP1912'),
    @('05:51:32', '07:11', '215A_EL PATO', 80, 'LP1912'),
    @('05:20:00', '07:12', '215A_EL PATO', 112, 'LP1912'),
    @('05:51:32', '07:15', '11_ETCHEVERRY', 84, 'LP1912'),
    @('06:52:23', '07:16', '16_SANTA ANA', 24, 'LP1912'),
    @('05:20:00', '07:16', '11_ETCHEVERRY', 116, 'LP1912'),
    @('05:51:32', '07:21', '26_HERNANDEZ', 90, 'LP1912'),
    @('06:23:52', '07:22', '10_OLMOS', 59, 'LP1912'),
    @('06:52:23', '07:23', '10_OLMOS', 31, 'LP1912'),
    @('06:52:23', '07:25', '10_OLMOS', 33, 'LP1912'),
    @('05:51:32', '07:28', '10_OLMOS', 97, 'LP1912'),
    @('05:51:32', '07:31', '11_ETCHEVERRY', 100, 'LP1912'),
    @('05:51:32', '07:31', '16_SANTA ANA', 100, 'LP1912'),
    @('05:51:32', '07:32', '84_COLONIA URQUIZA-ESC 49', 101, 'LP1912'),
    @('05:51:32', '07:36', '27_EL RETIRO', 105, 'LP1912'),
    @('06:52:23', '07:37', '27_EL RETIRO', 45, 'LP1912'),
    @('06:23:52', '07:38', '10_OLMOS', 75, 'LP1912'),
    @('05:51:32', '07:39', '10_OLMOS', 108, 'LP1912'),
    @('05:51:32', '07:47', '14_ABASTO', 116, 'LP1912'),
    @('07:46:15', '07:50', '14_ABASTO', 4, 'LP1912'),
    @('06:23:52', '07:51', '215D_EL PATO', 88, 'LP1912'),
    @('07:46:15', '07:55', '10_OLMOS', 9, 'LP1912'),
    @('07:46:15', '07:56', '16_SANTA ANA', 10, 'LP1912'),
    @('07:59:05', '08:02', '16_SANTA ANA', 3, 'LP1912'),
    @('06:52:23', '08:03', '23_HERNANDEZ', 71, 'LP1912'),
    @('06:23:52', '08:05', '23_HERNANDEZ', 102, 'LP1912'),
    @('07:46:15', '08:09', '11_ETCHEVERRY', 23, 'LP1912'),
    @('07:59:05', '08:11', '11_ETCHEVERRY', 12, 'LP1912'),
    @('06:23:52', '08:12', '15_ABASTO', 109, 'LP1912'),
    @('06:23:52', '08:20', '26_HERNANDEZ', 117, 'LP1912'),
    @('06:52:23', '08:21', '26_HERNANDEZ', 89, 'LP1912'),
    @('06:23:52', '08:22', '16_P MOR-SANTA ANA', 119, 'LP1912'),
    @('06:52:23', '08:23', '215B_EL PATO', 91, 'LP1912'),
    @('07:46:15', '08:23', '16_P MOR-SANTA ANA', 37, 'LP1912'),
    @('06:52:23', '08:27', '84_COLONIA URQUIZA-ESC 49', 95, 'LP1912'),
    @('07:46:15', '08:33', '10_OLMOS', 47, 'LP1912'),
    @('07:46:15', '08:33', '23_HERNANDEZ', 47, 'LP1912'),
    @('07:46:15', '08:34', '26_HERNANDEZ', 48, 'LP1912'),
    @('07:59:05', '08:35', '23_HERNANDEZ', 36, 'LP1912'),
    @('06:52:23', '08:42', '81_EL PELIGRO', 110, 'LP1912'),
    @('07:46:15', '08:44', '14_ABASTO', 58, 'LP1912'),
    @('07:59:05', '08:48', '26_HERNANDEZ', 49, 'LP1912'),
    @('07:59:05', '08:53', '10_OLMOS', 54, 'LP1912'),
    @('07:46:15', '08:54', '17_ROMERO', 68, 'LP1912'),
    @('07:46:15', '09:02', '215A_EL PATO', 76, 'LP1912'),
    @('07:46:15', '09:04', '11_ETCHEVERRY', 78, 'LP1912'),
    @('07:46:15', '09:11', '16_P MOR-SANTA ANA', 85, 'LP1912'),
    @('07:46:15', '09:17', '27_EL RETIRO', 91, 'LP1912'),
    @('07:46:15', '09:21', '26_HERNANDEZ', 95, 'LP1912'),
    @('07:46:15', '09:22', '16_SANTA ANA', 96, 'LP1912'),
    @('07:46:15', '09:23', '17_ROMERO', 97, 'LP1912'),
    @('07:59:05', '09:23', '16_SANTA ANA', 84, 'LP1912'),
    @('07:46:15', '09:24', '11_ETCHEVERRY', 98, 'LP1912'),
    @('07:46:15', '09:32', '15_ABASTO', 106, 'LP1912'),
    @('07:46:15', '09:33', '10_OLMOS', 107, 'LP1912'),
    @('07:46:15', '09:42', '215C_EL PATO', 116, 'LP1912'),
    @('07:59:05', '09:44', '14_ABASTO', 105, 'LP1912'),
    @('07:59:05', '09:52', '15_ABASTO', 113, 'LP1912')
)

for ($i = 0; $i -lt $data1.Length; $i++) {
    $r = $i + 6
    $row = $data1[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------- Sheet: LP1912-215 ----------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 07:59:05'

# ---------- Sheet: 6203-6173 ----------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 07:59:05'
$ws3.Range("A3").Value = 'Total filas: 22'

$data3 = @(
    @('04:36:47', '05:43', '215A_LA PLATA', 67, 'L6173'),
    @('04:01:01', '05:44', '215A_LA PLATA', 103, 'L6173'),
    @('04:51:28', '06:09', '215A_LA PLATA', 78, 'L6173'),
    @('04:36:47', '06:10', '215A_LA PLATA', 94, 'L6173'),
    @('04:36:47', '06:32', '215C_LA PLATA', 116, 'L6203'),
    @('04:51:28', '06:33', '215C_LA PLATA', 102, 'L6203'),
    @('06:23:52', '06:59', '215B_LP-P MOR-1 Y 57', 36, 'L6173'),
    @('05:20:00', '07:00', '215B_LP-P MOR-1 Y 57', 100, 'L6173'),
    @('05:51:32', '07:35', '215A_LA PLATA', 104, 'L6173'),
    @('06:52:23', '07:38', '215A_LA PLATA', 46, 'L6173'),
    @('07:16:53', '07:44', '215A_LA PLATA', 28, 'L6173'),
    @('07:46:15', '07:51', '215A_LA PLATA', 5, 'L6173'),
    @('06:23:52', '08:06', '215C_LA PLATA', 103, 'L6203'),
    @('07:46:15', '08:09', '215C_LA PLATA', 23, 'L6203'),
    @('07:16:53', '08:10', '215C_LA PLATA', 54, 'L6203'),
    @('06:52:23', '08:11', '215C_LA PLATA', 79, 'L6203'),
    @('07:59:05', '08:15', '215C_LA PLATA', 16, 'L6203'),
    @('06:52:23', '08:40', '215A_LA PLATA', 108, 'L6173'),
    @('07:46:15', '08:45', '215A_LA PLATA', 59, 'L6173'),
    @('07:59:05', '08:46', '215A_LA PLATA', 47, 'L6173'),
    @('07:16:53', '09:08', '215D_LA PLATA', 112, 'L6203'),
    @('07:46:15', '09:09', '215D_LA PLATA', 83, 'L6203')
)

for ($i = 0; $i -lt $data3.Length; $i++) {
    $r = $i + 6
    $row = $data3[$i]
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
}

Write-Output "Update complete"